# 5311: cleaned up functions more, continued merging with loc
#
# Fixes the "Klamath Trinity Non-Emergency Transportation" row in the
# bad_fuzzy_matches crosswalk:
#   - row 17 (Trinity County Dept of Transportation) picks up the
#     ntd_id "9R02-91057" that had been mis-attached elsewhere, and its
#     stray header-style shading is cleared
#   - row 21's org name is corrected from a mojibake'd zero-width space
#     to a real one, and it is no longer merged with an agency/ntd_id
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cell value fixes -------------------------------------------------

# Row 17 gains the ntd_id that belongs to the Trinity County dept row.
$ws.Range("C17").Value = "9R02-91057"

# Row 21's organization name: replace the mis-encoded zero-width space
# with a genuine U+200B so the text renders clean.
$ws.Range("A21").Value = "Klamath Trinity Non-Emergency Transportation" + [char]0x200B

# --- formatting fixes ---------------------------------------------------

# Row 17 no longer carries the section-header shading; keep its font
# size (12pt, matching the rest of the data rows) but drop the fill.
$hdrRow = $ws.Range("A17:C17")
$hdrRow.Interior.Pattern = -4142   # xlPatternNone
$hdrRow.Font.Size = 12

# --- selection -----------------------------------------------------------

$ws.Range("D24").Select() | Out-Null
